$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Rename countries to match updated naming conventions
$ws.Range("A2").Value = "Russian Federation"
$ws.Range("A15").Value = "Czechia"

# Restore the active selection to A5
$ws.Range("A5").Select()
